# Generate Report for handoff
# - Overview: "Ready for handoff" -> "Handoff transform failed"
# - zh-cn / de-de sheets: the handoff transform failed, so the per-language
#   rows no longer have a produced handoff file / timestamp, and the
#   dependency status flips from "Include" to "Ignored".

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handoff transform failed"
$overview.Range("C2").Value = "Handoff transform failed"

$langSheets = @("zh-cn", "de-de")

foreach ($sheetName in $langSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Status column (B) mirrors the Overview status text.
    $ws.Range("B2").Value = "Handoff transform failed"

    # The handoff transform failed, so there is no Latest Handoff File
    # anymore for row 2 - drop its value and its hyperlink.
    $ws.Range("C2").Value = ""

    # Latest Handoff Datetime reverts to the "never happened" sentinel.
    $ws.Range("D2").Value = "0001-01-01 00:00:00"
    $ws.Range("D3").Value = "0001-01-01 00:00:00"

    # Latest Handback DateTime sentinel stays the same on both rows.
    $ws.Range("G2").Value = "0001-01-01 00:00:00"
    $ws.Range("G3").Value = "0001-01-01 00:00:00"

    # Handoff Reason flips from Include/Ignored to Ignored/Ignored.
    $ws.Range("H2").Value = "Ignored"
    $ws.Range("H3").Value = "Ignored"

    # Rebuild the hyperlinks collection: only A2 (source md file) and A3
    # (.localization-config) keep links once the handoff-file link is gone.
    $mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/f341b47a67180c1da61c4e4c4a25d3b69af74dbd/e2e/a7121821-1a0d-4dfc-9a39-c640b6860ad5.md"
    $configAddress = "https://github.com/OpenLocalizationTest/oltest/blob/f341b47a67180c1da61c4e4c4a25d3b69af74dbd/.localization-config"

    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $mdAddress, [Type]::Missing, [Type]::Missing, "a7121821-1a0d-4dfc-9a39-c640b6860ad5.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), $configAddress, [Type]::Missing, [Type]::Missing, ".localization-config")
}
